# Update symbol list values in cryptos worksheet (GitHub Actions refresh).
# All target cells are stored as text (inlineStr) in the workbook, so for the
# "Price" column (D) -- whose new values look numeric -- we force the cell's
# NumberFormat to Text ("@") before assigning the value. This prevents Excel
# from silently re-interpreting strings like "0.8140" or "21.90" as numbers
# (which would drop significant trailing zeros / change precision).
# Columns B, C and E contain non-numeric text already, so a plain assignment
# keeps them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.74"

# Row 3 - OKB
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.92"

# Row 4 - HuobiToken
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.386"

# Row 5 - Cronos
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05991"

# Row 7 - MXToken
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8140"

# Row 8 - FTXToken
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9563"

# Row 9 - WazirX
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1426"

# Row 10 - MandalaExchangeToken
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07434"

# Row 11 - LiechtensteinCryptoassetsExchange
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03334"

# Row 12 - BitrueCoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03055"

# Row 13 - BitMartToken
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09417"

# Row 15 - BitForexToken
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001587"

# Row 16 - CoinExToken
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04797"

# Row 17 - One
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005902"
$ws.Range("E17").Value = "16OneONEWorstin24h"

# Row 18 - TigerCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006122"

# Row 20 - BitKan
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009870"

# Row 22 - LEO
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.686"

# Row 23 - KuCoinToken
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.433"

# Row 24 - BTSEToken
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"

# Row 40 - IDEX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03991"

# Row 41 - was BKEXToken, now KickToken (rows 41-43 rotate)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006532"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 - was CEJI, now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 - was KickToken, now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005777"

# Row 45 - CoinLion
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005278"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8503"

# Row 48 - BOLO
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01260"
$ws.Range("E48").Value = "47BOLOBOLO"
